# edit.ps1 - Apply commit "Enhance: Preliminary discussion should now be complete,
# now onto practical procedures for CAS calc usage." to the VCE CAS transcript.
#
# Real (visible) content changes only - the diff's many proofErr / gramStart / gramEnd
# / spellStart / spellEnd removals and run re-splits carry no visible text change, so
# they are left to Word's own bookkeeping and are not separately reproduced here.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "... VCAA has already approved at least four software applications ..."
#    -> "... VCAA has already approved around half a dozen software applications ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "VCAA has already approved at least four software applications",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "VCAA has already approved around half a dozen software applications", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "... in order to helps us do CAS Maths, or Mathematical Methods in years 11 ..."
#    -> "... in order to help us do CAS Maths, or what is officially called
#         Mathematical Methods in years 11 ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "in order to helps us do CAS Maths, or Mathematical Methods in years 11",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in order to help us do CAS Maths, or what is officially called Mathematical Methods in years 11", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append a closing sentence to the "Using Software Calculators" paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "perform calculations from examples or exercises in a textbook or workbook.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "perform calculations from examples or exercises in a textbook or workbook.  Hopefully these insights get you off to a good start.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Remove the two "Nathaniel's reflections ..." bullet paragraphs that used
#    to sit right before the "Available and approved CAS software calculators:"
#    list (found by locating that anchor paragraph and walking back up).
# ---------------------------------------------------------------------------
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Available and approved CAS software calculators:*") {
        $anchor = $i
        break
    }
}
if ($anchor -ne $null) {
    # Delete the two paragraphs immediately preceding the anchor paragraph.
    $d.Paragraphs.Item($anchor - 2).Range.Delete() | Out-Null
    $d.Paragraphs.Item($anchor - 2).Range.Delete() | Out-Null
}

# ---------------------------------------------------------------------------
# 5) Fix the "Flexability" typo -> "Flexibility".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Flexability", $true, $false, $false, $false, $false, $true, 1, $false,
    "Flexibility", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Strip the "•<tab>" / "o<tab>" outline markers from every remaining list
#    paragraph - the list is no longer manually bulleted.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.Length -ge 2 -and $t[1] -eq [char]9 -and ($t[0] -eq [char]8226 -or $t[0] -eq 'o')) {
        $r = $p.Range
        $d.Range($r.Start, $r.Start + 2).Delete() | Out-Null
    }
}

# ---------------------------------------------------------------------------
# 7) Insert the new "Installation files / process." list item right after
#    "Compatibility across common devices and platforms."
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "Compatibility across common devices and platforms.") {
        $d.Paragraphs.Item($i).Range.InsertParagraphAfter() | Out-Null
        $d.Paragraphs.Item($i + 1).Range.Text = "Installation files / process."
        break
    }
}

$d.Save()
